# ---------------------------------------------------------------------------
# Switch the presentation's design theme from "Integral" (Red Violet colour
# scheme) to the built-in "Office Theme" (Office colour scheme), i.e. the
# Design gallery action a user would perform from PowerPoint's Design tab.
#
# The deck has a single slide master, so its theme part (ppt/theme/theme1.xml)
# is shared by every slide; updating the 12 theme-colour slots through any
# slide's ThemeColorScheme updates that one shared theme part.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme (RGB() values == little-endian hex of the
# target sRGB colour, e.g. 6968388 = &H44546A == dk2 "44546A"):
$tcs.Colors(1).RGB  = 0          # Dark 1    (dk1)      000000
$tcs.Colors(2).RGB  = 16777215   # Light 1   (lt1)      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # Dark 2    (dk2)      44546A
$tcs.Colors(4).RGB  = 15132391   # Light 2   (lt2)      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # Accent 1  (accent1)  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # Accent 2  (accent2)  ED7D31
$tcs.Colors(7).RGB  = 10855845   # Accent 3  (accent3)  A5A5A5
$tcs.Colors(8).RGB  = 49407      # Accent 4  (accent4)  FFC000
$tcs.Colors(9).RGB  = 12874308   # Accent 5  (accent5)  4472C4
$tcs.Colors(10).RGB = 4697456    # Accent 6  (accent6)  70AD47
$tcs.Colors(11).RGB = 12673797   # Hyperlink (hlink)    0563C1
$tcs.Colors(12).RGB = 7491477    # Followed Hyperlink (folHlink) 954F72
